$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "ECs" sending-cluster block (original rows 2-4).
# The remaining rows (originally 5-10, the FAPs/MuSCs blocks) shift up to rows 2-7.
$ws.Range("A2:T4").EntireRow.Delete()

# New TPM-derived values for the remaining rows (now rows 2-7).
$data = @(
    @("FAPs",  "Col2a1", "Itga10", "ECs",   1, 0.3333333333333333, 0.2021233333333333, 0.60637,    0.9764223557676824,  0.9764223557676823,  3, 1, 0.353476,            1.060428, 0.07277959798226569, 0.07277959798226567, 0.07144574737333333, 0.6430117263599999, 0.07106362651366872, 0.07106362651366871),
    @("FAPs",  "Col2a1", "Itga10", "FAPs",  1, 0.3333333333333333, 0.2021233333333333, 0.60637,    0.9764223557676824,  0.9764223557676823,  3, 1, 1.358024333333333,   4.074073, 0.2796129441040817,  0.2796129441040817,  0.2744884050011111,  2.47039564501,      0.2730203295852448,  0.2730203295852447),
    @("FAPs",  "Col2a1", "Itga10", "MuSCs", 1, 0.3333333333333333, 0.2021233333333333, 0.60637,    0.9764223557676824,  0.9764223557676823,  3, 1, 3.145300333333334,   9.435901000000001, 0.6476074579136527, 0.6476074579136526, 0.6357385877077778, 5.721647289370001, 0.632338399668769,  0.6323383996687688),
    @("MuSCs", "Col2a1", "Itga10", "ECs",   2, 0.6666666666666666, 0.004880666666666667, 0.014642, 0.02357764423231757, 0.02357764423231757, 3, 1, 0.353476,            1.060428, 0.07277959798226569, 0.07277959798226567, 0.001725198530666666, 0.015526786776,    0.001715971468596958, 0.001715971468596958),
    @("MuSCs", "Col2a1", "Itga10", "FAPs",  2, 0.6666666666666666, 0.004880666666666667, 0.014642, 0.02357764423231757, 0.02357764423231757, 3, 1, 1.358024333333333,   4.074073, 0.2796129441040817,  0.2796129441040817,  0.006628064096222223, 0.05965257686600001, 0.006592614518836938, 0.006592614518836937),
    @("MuSCs", "Col2a1", "Itga10", "MuSCs", 2, 0.6666666666666666, 0.004880666666666667, 0.014642, 0.02357764423231757, 0.02357764423231757, 3, 1, 3.145300333333334,   9.435901000000001, 0.6476074579136527, 0.6476074579136526, 0.01535116249355556, 0.138160462442,      0.01526905824488368, 0.01526905824488368)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $colIndex = $j + 1
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$j]
    }
}
